$d = $word.ActiveDocument

# 1. Update the name/title line to add "- Software Engineer"
$d.Content.Find.Execute("Darrel Daquigan ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Darrel Daquigan - Software Engineer ", 2)

# 2. Update the objective line (drop "/web" and "or development")
$d.Content.Find.Execute( `
    "Looking for an entry-level software/web engineering or development position to utilize my degree in computer science", `
    $true, $false, $false, $false, $false, `
    $true, 1, $false, `
    "Looking for an entry-level software engineering position to utilize my degree in computer science", 2)

# 3. Add Python to the skills/experience list
$d.Content.Find.Execute("jQuery, MATLAB, Prolog, Scheme/Racket, Ruby, R", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "jQuery, Python, MATLAB, Prolog, Scheme/Racket, Ruby, R", 2)

# 4. Add "Algorithm analysis" right after " Numerical analysis," in the Other Skills list.
#    The comma that immediately follows " Numerical analysis" sits in its own run, so
#    locate it precisely (rather than a blind text Find/Replace) and rewrite just that
#    run's text, matching the source edit exactly.
$r2 = $d.Content
$r2.Find.Execute(" Numerical analysis", $true, $false, $false, $false, $false, `
                  $true, 1, $false, "", 0)
$commaRange = $d.Range($r2.End, $r2.End + 1)
$commaRange.Text = ", Algorithm analysis, "

# 5. Fix quotation marks around X language (curly double quotes -> curly single quotes)
$d.Content.Find.Execute( `
    "debugger for " + [char]0x201C + "X" + [char]0x201D + " language", `
    $true, $false, $false, $false, $false, `
    $true, 1, $false, `
    "debugger for " + [char]0x2018 + "X" + [char]0x2019 + " language", 2)
